$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("9_10", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
$r.Text = "10"
$r.Collapse(0)
$r.InsertAfter("_1")
$r.Collapse(0)
$r.InsertAfter("1")
